$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.198.82'
$ws.Range("E2").Value = '  -0.67%  '
$ws.Range("D3").Value = '1.826.16'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.42'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5976'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.51%  '
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07009'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.76%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2776'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.13%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.37'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.46%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07646'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.97%  '
$ws.Range("D12").Value = '1.839.28'
$ws.Range("E12").Value = '  -0.03%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.778'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.03%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.000009922'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.51%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6233'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -7.78%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '78.67'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.95%  '
$ws.Range("D17").Value = '29.187.84'
$ws.Range("E17").Value = '  -0.66%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.817'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -6.62%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '222.92'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.41%  '
$ws.Range("E20").Value = '  +0.18%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.63'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.60%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.982'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.78%  '
$ws.Range("E23").Value = '  +0.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '154.73'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.946'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.37%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1290'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.53%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.48'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.93%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.478'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.68%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.06194'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -13.78%  '
$ws.Range("E30").Value = '  -2.90%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.826'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.781'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.54%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.104'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.13%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.731'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.02%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6429'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.86%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.541'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = '1.220.37'
$ws.Range("E37").Value = '  -1.14%  '
$ws.Range("E38").Value = '  -3.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.512'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.08%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01725'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.27%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8977'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.28%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.002'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.17%  '
$ws.Range("D43").Value = '1.980.01'
$ws.Range("E43").Value = '  -1.47%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.02'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.98%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '62.19'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.06%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000116'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.05%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.482'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4555'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.60%  '
$ws.Range("E49").Value = '  -9.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05492'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.96%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.389'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -8.30%  '
